$d = $word.ActiveDocument

# --- Change 1: append a new sentence/run to the "Potential solutions..." paragraph ---
$pPotential = $d.Paragraphs(34)
$rPotential = $pPotential.Range
$rPotential.End = $rPotential.End - 1   # exclude the paragraph mark
$rPotential.Collapse(0)
$rPotential.InsertAfter("Also, examining the pattern and figuring it out and how it repeats.")

# --- Change 2: remove the _GoBack bookmark from the last paragraph, add two new paragraphs,
#     and re-create the bookmark at the end of the new final paragraph ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$pLast = $d.Paragraphs($d.Paragraphs.Count)
$rLast = $pLast.Range
$rLast.Collapse(0)
$rLast.InsertParagraphAfter()

$pEmpty = $d.Paragraphs($d.Paragraphs.Count - 1)
$rEmpty = $pEmpty.Range
$rEmpty.Collapse(0)
$rEmpty.InsertParagraphAfter()

$pNew = $d.Paragraphs($d.Paragraphs.Count)
$rNew = $pNew.Range
$rNew.Collapse(0)
$rNew.InsertAfter([char]9 + "The solution I would use is examining the pattern. Because the pattern is 5 fingers in a 10 count the pattern remains the same, regardless of how many times it is repeated. This would work in all cases of this problem, provided the amount of fingers counted and the count itself remained the same. ")

$pFinal = $d.Paragraphs($d.Paragraphs.Count)
$rFinal = $pFinal.Range
$rFinal.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rFinal)
